# Add a new row (ID=2) describing "Purple bloom", and trim the trailing
# "Hi Sara" from the existing Endless Wind description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Row 3 data (set title/description text first so the shared-string
# table is built up in the same order as the source workbook).
$ws.Range("B3").Value = "Purple bloom"
$ws.Range("E3").Value = "Sun rising beneath an endless purple wind."

# Update Row 2's description: drop the trailing "Hi Sara".
$ws.Range("E2").Value = "Beautiful interactions of strings, weaving endlessly into the infinite."

$ws.Range("A3").Value = 2
$ws.Range("C3").Value = "Jad Okaily"
$ws.Range("D3").Value = 2025

# Move the active selection down to E4, matching the new sheet state.
$ws.Range("E4").Select() | Out-Null
